$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original cell style, force text format so numeric-looking
# price strings are not auto-converted to numbers by Excel, then restore
# the original style so no formatting change is introduced.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.007.10"
Set-TextValue $ws.Range("E2") "  -0.81%  "

Set-TextValue $ws.Range("D3") "1.903.68"
Set-TextValue $ws.Range("E3") "  -0.61%  "

Set-TextValue $ws.Range("D4") "0.9997"
Set-TextValue $ws.Range("E4") "  -0.36%  "

Set-TextValue $ws.Range("D5") "0.7478"
Set-TextValue $ws.Range("E5") "  +0.52%  "

Set-TextValue $ws.Range("D6") "241.50"
Set-TextValue $ws.Range("E6") "  -0.48%  "

Set-TextValue $ws.Range("D7") "0.9999"
Set-TextValue $ws.Range("E7") "  -0.41%  "

Set-TextValue $ws.Range("D8") "0.3076"
Set-TextValue $ws.Range("E8") "  -1.68%  "

Set-TextValue $ws.Range("D9") "25.61"
Set-TextValue $ws.Range("E9") "  -5.81%  "

Set-TextValue $ws.Range("D10") "0.06898"
Set-TextValue $ws.Range("E10") "  -0.61%  "

Set-TextValue $ws.Range("D11") "0.08015"
Set-TextValue $ws.Range("E11") "  +0.24%  "

Set-TextValue $ws.Range("D12") "0.7559"
Set-TextValue $ws.Range("E12") "  -1.62%  "

Set-TextValue $ws.Range("D13") "1.899.84"
Set-TextValue $ws.Range("E13") "  -1.12%  "

Set-TextValue $ws.Range("D14") "5.271"
Set-TextValue $ws.Range("E14") "  -0.29%  "

Set-TextValue $ws.Range("D15") "91.57"
Set-TextValue $ws.Range("E15") "  +0.41%  "

Set-TextValue $ws.Range("D16") "6.178"
Set-TextValue $ws.Range("E16") "  +5.95%  "

Set-TextValue $ws.Range("D17") "30.003.21"
Set-TextValue $ws.Range("E17") "  -0.97%  "

Set-TextValue $ws.Range("D18") "14.05"
Set-TextValue $ws.Range("E18") "  -0.76%  "

Set-TextValue $ws.Range("D19") "0.000007757"
Set-TextValue $ws.Range("E19") "  -1.05%  "

Set-TextValue $ws.Range("D20") "237.53"
Set-TextValue $ws.Range("E20") "  -3.75%  "

Set-TextValue $ws.Range("D21") "0.9998"
Set-TextValue $ws.Range("E21") "  -0.47%  "

Set-TextValue $ws.Range("D22") "2.151.37"
Set-TextValue $ws.Range("E22") "  -0.48%  "

Set-TextValue $ws.Range("D23") "0.9997"
Set-TextValue $ws.Range("E23") "  -0.24%  "

Set-TextValue $ws.Range("D24") "7.092"
Set-TextValue $ws.Range("E24") "  +7.60%  "

Set-TextValue $ws.Range("D25") "9.324"
Set-TextValue $ws.Range("E25") "  -0.49%  "

Set-TextValue $ws.Range("D26") "166.39"
Set-TextValue $ws.Range("E26") "  +0.80%  "

Set-TextValue $ws.Range("D27") "18.82"
Set-TextValue $ws.Range("E27") "  -0.31%  "

Set-TextValue $ws.Range("D28") "0.1266"
Set-TextValue $ws.Range("E28") "  -1.01%  "

Set-TextValue $ws.Range("D29") "2.063"
Set-TextValue $ws.Range("E29") "  -3.68%  "

Set-TextValue $ws.Range("E30") "  -1.72%  "

Set-TextValue $ws.Range("D31") "1.527"
Set-TextValue $ws.Range("E31") "  -1.15%  "

Set-TextValue $ws.Range("D32") "4.310"
Set-TextValue $ws.Range("E32") "  -0.56%  "

Set-TextValue $ws.Range("D33") "4.051"
Set-TextValue $ws.Range("E33") "  -0.18%  "

Set-TextValue $ws.Range("D34") "0.05416"
Set-TextValue $ws.Range("E34") "  +4.75%  "

Set-TextValue $ws.Range("D35") "1.288"
Set-TextValue $ws.Range("E35") "  -0.75%  "

Set-TextValue $ws.Range("D36") "0.7385"
Set-TextValue $ws.Range("E36") "  -0.76%  "

Set-TextValue $ws.Range("E37") "  -1.57%  "

Set-TextValue $ws.Range("D38") "0.01945"
Set-TextValue $ws.Range("E38") "  +0.62%  "

Set-TextValue $ws.Range("D39") "2.758"
Set-TextValue $ws.Range("E39") "  -0.41%  "

Set-TextValue $ws.Range("D40") "6.242"
Set-TextValue $ws.Range("E40") "  -2.50%  "

Set-TextValue $ws.Range("D41") "0.4459"
Set-TextValue $ws.Range("E41") "  +0.17%  "

Set-TextValue $ws.Range("D42") "72.78"
Set-TextValue $ws.Range("E42") "  -3.87%  "

Set-TextValue $ws.Range("D43") "1.948"
Set-TextValue $ws.Range("E43") "  +0.73%  "

Set-TextValue $ws.Range("D44") "1.000"
Set-TextValue $ws.Range("E44") "  -0.33%  "

Set-TextValue $ws.Range("D45") "0.8314"
Set-TextValue $ws.Range("E45") "  -0.67%  "

Set-TextValue $ws.Range("D46") "7.701"
Set-TextValue $ws.Range("E46") "  +1.72%  "

Set-TextValue $ws.Range("D47") "101.61"
Set-TextValue $ws.Range("E47") "  +0.41%  "

Set-TextValue $ws.Range("D48") "9.864"
Set-TextValue $ws.Range("E48") "  +0.88%  "

Set-TextValue $ws.Range("D49") "2.057.53"
Set-TextValue $ws.Range("E49") "  -0.55%  "

Set-TextValue $ws.Range("D50") "36.58"
Set-TextValue $ws.Range("E50") "  -1.28%  "

Set-TextValue $ws.Range("D51") "0.1166"
Set-TextValue $ws.Range("E51") "  -3.54%  "
